$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add headers I1 and J1, copying header style/format from G1
$ws.Range("G1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Populate I0 (col I) and IF (col J) data for rows 2-74
$ijData = @{
    2 = @(6, 6)
    3 = @(8, 8)
    4 = @(7, 7)
    5 = @(7, 7)
    6 = @(8, 8)
    7 = @(7, 7)
    8 = @(7, 7)
    9 = @(8, 8)
    10 = @(7, 7)
    11 = @(8, 8)
    12 = @(8, 8)
    13 = @(9, 9)
    14 = @(8, 8)
    15 = @(8, 8)
    16 = @(9, 9)
    17 = @(7, 7)
    18 = @(8, 8)
    19 = @(9, 9)
    20 = @(7, 7)
    21 = @(7, 7)
    22 = @(8, 8)
    23 = @(9, 9)
    24 = @(8, 8)
    25 = @(9, 9)
    26 = @(9, 8)
    27 = @(9, 9)
    28 = @(8, 8)
    29 = @(8, 8)
    30 = @(8, 8)
    31 = @(9, 9)
    32 = @(9, 9)
    33 = @(9, 9)
    34 = @(9, 9)
    35 = @(9, 9)
    36 = @(9, 9)
    37 = @(9, 9)
    38 = @(9, 9)
    39 = @(8, 8)
    40 = @(9, 9)
    41 = @(9, 9)
    42 = @(10, 10)
    43 = @(9, 9)
    44 = @(9, 9)
    45 = @(9, 9)
    46 = @(9, 9)
    47 = @(9, 9)
    48 = @(9, 9)
    49 = @(8, 9)
    50 = @(9, 9)
    51 = @(9, 9)
    52 = @(9, 9)
    53 = @(9, 9)
    54 = @(9, 10)
    55 = @(9, 9)
    56 = @(9, 9)
    57 = @(9, 9)
    58 = @(9, 9)
    59 = @(9, 9)
    60 = @(9, 9)
    61 = @(9, 9)
    62 = @(9, 9)
    63 = @(9, 9)
    64 = @(9, 9)
    65 = @(9, 9)
    66 = @(9, 9)
    67 = @(9, 9)
    68 = @(9, 9)
    69 = @(9, 9)
    70 = @(9, 9)
    71 = @(6, 6)
    72 = @(6, 6)
    73 = @(4, 4)
    74 = @(3, 3)
}

foreach ($row in $ijData.Keys) {
    $vals = $ijData[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
